$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Bug #14 ("Menu only displays 1 button instad of 3") is now fixed.
#    Row 18 stops being the last row in the table, so it picks up the same
#    "interior" double-border look that every other fixed-bug row (4-17) has,
#    and its "Date Fixed" cell (L18) gets filled in.
# ---------------------------------------------------------------------------

# Copy the formatting of the row above (still a normal interior row) onto
# row 18 - this flips its borders/row-thickness from "last row" to
# "interior row" (double top+bottom separators) without touching values.
$ws.Range("A17:L17").Copy()
$ws.Range("A18:L18").PasteSpecial(-4122)

# Row 18 was fixed the same day it was found.
$ws.Range("L18").Value = 42484

# ---------------------------------------------------------------------------
# 2) A new bug (#15) was found the same day - it goes into row 19, which was
#    just an empty template row. Give it the same look the old "last row"
#    (old row 18) used to have.
# ---------------------------------------------------------------------------

$ws.Range("A19:L19").ClearContents()
$ws.Range("A18:L18").Copy()
$ws.Range("A19:L19").PasteSpecial(-4122)

$ws.Range("A19").Value = 15
$ws.Range("B19").Value = 42484
$ws.Range("C19").Value = "Charles"
$ws.Range("D19").Value = "Code - Function"
$ws.Range("E19").Value = "High"
$ws.Range("F19").Value = "Charles"
$ws.Range("G19").Value = "Clicking play from the menu starts the game, but the level acts as if it started playing when the menu started"
$ws.Range("L19").Value = ""

# ---------------------------------------------------------------------------
# 3) Scroll / selection bookkeeping to match where the author was working.
# ---------------------------------------------------------------------------

$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("A19:L19").Select()
